# Scheduled market-data refresh: update cached currentAveragePrice /
# LevePrice / LeveProfit columns (H:N) for the affected Leve rows.
# Values come from an external price snapshot; no formulas are used
# in this workbook, so cells are written as literals.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 7: The Bleat Is On
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 14: Wand-full Tonight
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
# Row 64: Forged from the Void
$ws.Range("H64").Value = 3367.125
$ws.Range("I64").Value = 3491
$ws.Range("J64").Value = 2500
$ws.Range("K64").Value = 3491
$ws.Range("L64").Value = 2500
$ws.Range("M64").Value = -3243
$ws.Range("N64").Value = -2996
# Row 67: Dodging the Draft (L)
$ws.Range("H67").Value = 3367.125
$ws.Range("I67").Value = 3491
$ws.Range("J67").Value = 2500
$ws.Range("K67").Value = 3491
$ws.Range("L67").Value = 2500
$ws.Range("M67").Value = -2633
$ws.Range("N67").Value = -4216
# Row 98: The Dotted Line
$ws.Range("H98").Value = 28461842
$ws.Range("I98").Value = 13334920
$ws.Range("J98").Value = 51152228
$ws.Range("K98").Value = 13334920
$ws.Range("L98").Value = 51152228
$ws.Range("M98").Value = -13333422
$ws.Range("N98").Value = -51155224
# Row 122: Wishful Inking
$ws.Range("H122").Value = 28461842
$ws.Range("I122").Value = 13334920
$ws.Range("J122").Value = 51152228
$ws.Range("K122").Value = 40004760
$ws.Range("L122").Value = 153456684
$ws.Range("M122").Value = -40002310
$ws.Range("N122").Value = -153461584
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 2416309.2
$ws.Range("I132").Value = 868.3415
$ws.Range("K132").Value = 2605.0245
$ws.Range("M132").Value = -75.02449999999999
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 4032.375
$ws.Range("I141").Value = 4032.375
$ws.Range("K141").Value = 12097.125
$ws.Range("M141").Value = -6917.125

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 44675
$ws.Range("I2").Value = 50914.285
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 50914.285
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -50801.285
$ws.Range("N2").Value = -1226
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 8229549
$ws.Range("I32").Value = 2124541.5
$ws.Range("K32").Value = 2124541.5
$ws.Range("M32").Value = -2124254.5
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4020466.2
$ws.Range("I61").Value = 1737005.2
$ws.Range("K61").Value = 1737005.2
$ws.Range("M61").Value = -1736793.2
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 77384260
$ws.Range("I74").Value = 84524950
$ws.Range("J74").Value = 66673212
$ws.Range("K74").Value = 84524950
$ws.Range("L74").Value = 66673212
$ws.Range("M74").Value = -84524076
$ws.Range("N74").Value = -66674960
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 77384260
$ws.Range("I77").Value = 84524950
$ws.Range("J77").Value = 66673212
$ws.Range("K77").Value = 422624750
$ws.Range("L77").Value = 333366060
$ws.Range("M77").Value = -422620382
$ws.Range("N77").Value = -333374796
# Row 116: No Scope
$ws.Range("H116").Value = 44675
$ws.Range("I116").Value = 50914.285
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 50914.285
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = -48620.285
$ws.Range("N116").Value = -5588
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 12157312
$ws.Range("I132").Value = 12350189
$ws.Range("K132").Value = 37050567
$ws.Range("M132").Value = -37048037
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4020466.2
$ws.Range("I136").Value = 1737005.2
$ws.Range("K136").Value = 5211015.6
$ws.Range("M136").Value = -5208465.6

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 44675
$ws.Range("I3").Value = 50914.285
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 50914.285
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -50800.285
$ws.Range("N3").Value = -1228
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 20015558
$ws.Range("I20").Value = 41684420
$ws.Range("J20").Value = 13529.154
$ws.Range("K20").Value = 41684420
$ws.Range("L20").Value = 13529.154
$ws.Range("M20").Value = -41684173
$ws.Range("N20").Value = -14023.154

$ws = $wb.Worksheets.Item("CRP")
# Row 107: Built to Last
$ws.Range("H107").Value = 566.4091
$ws.Range("I107").Value = 280.08334
$ws.Range("J107").Value = 910
$ws.Range("K107").Value = 280.08334
$ws.Range("L107").Value = 910
$ws.Range("M107").Value = 1639.91666
$ws.Range("N107").Value = -4750

$ws = $wb.Worksheets.Item("CUL")
# Row 104: Fits to a Tea
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("M104").ClearContents()
# Row 129: Comfort Food
$ws.Range("H129").Value = 1856.4242
$ws.Range("I129").Value = 1576.25
$ws.Range("J129").Value = 1946.08
$ws.Range("K129").Value = 4728.75
$ws.Range("L129").Value = 5838.24
$ws.Range("M129").Value = 271.25
$ws.Range("N129").Value = -15838.24

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 6947093.5
$ws.Range("I122").Value = 2840.2354
$ws.Range("K122").Value = 8520.706200000001
$ws.Range("M122").Value = -6070.706200000001

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 831.6667
$ws.Range("I7").Value = 818
$ws.Range("K7").Value = 818
$ws.Range("M7").Value = -706
# Row 19: Targe Up
$ws.Range("H19").Value = 445
$ws.Range("I19").Value = 445
$ws.Range("K19").Value = 445
$ws.Range("M19").Value = -275
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 4180.846
$ws.Range("I22").Value = 567.5
$ws.Range("J22").Value = 9962.200000000001
$ws.Range("K22").Value = 567.5
$ws.Range("L22").Value = 9962.200000000001
$ws.Range("M22").Value = -272.5
$ws.Range("N22").Value = -10552.2
# Row 27: Fire and Hide
$ws.Range("H27").Value = 4180.846
$ws.Range("I27").Value = 567.5
$ws.Range("J27").Value = 9962.200000000001
$ws.Range("K27").Value = 567.5
$ws.Range("L27").Value = 9962.200000000001
$ws.Range("M27").Value = -460.5
$ws.Range("N27").Value = -10176.2
# Row 126: Battered Books
$ws.Range("H126").Value = 831.6667
$ws.Range("I126").Value = 818
$ws.Range("K126").Value = 2454
$ws.Range("M126").Value = 16

$ws = $wb.Worksheets.Item("WVR")
# Row 64: Ribbon of Remembrance
$ws.Range("H64").Value = 23250
$ws.Range("J64").Value = 23250
$ws.Range("L64").Value = 23250
$ws.Range("N64").Value = -23746
# Row 67: The Road Was a Ribbon of Moonlight (L)
$ws.Range("H67").Value = 23250
$ws.Range("J67").Value = 23250
$ws.Range("L67").Value = 23250
$ws.Range("N67").Value = -24966
# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 22514.285
$ws.Range("I81").Value = 4000
$ws.Range("J81").Value = 23938.46
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 47876.92
$ws.Range("M81").Value = -6939
$ws.Range("N81").Value = -49998.92
# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 22514.285
$ws.Range("I84").Value = 4000
$ws.Range("J84").Value = 23938.46
$ws.Range("K84").Value = 40000
$ws.Range("L84").Value = 239384.6
$ws.Range("M84").Value = -34696
$ws.Range("N84").Value = -249992.6
